$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 656.6667
$ws.Range("I2").Value = 656.6667
$ws.Range("K2").Value = 656.6667
$ws.Range("M2").Value = -543.6667
$ws.Range("H43").Value = 2037.75
$ws.Range("I43").Value = 2250.5
$ws.Range("J43").Value = 1966.8334
$ws.Range("K43").Value = 2250.5
$ws.Range("L43").Value = 1966.8334
$ws.Range("M43").Value = -2181.5
$ws.Range("N43").Value = -2104.8334
$ws.Range("H64").Value = 3060.375
$ws.Range("J64").Value = 3085.111
$ws.Range("L64").Value = 3085.111
$ws.Range("N64").Value = -3581.111
$ws.Range("H67").Value = 3060.375
$ws.Range("J67").Value = 3085.111
$ws.Range("L67").Value = 3085.111
$ws.Range("N67").Value = -4801.111
$ws.Range("H111").Value = 2415.7144
$ws.Range("I111").Value = 2415.6
$ws.Range("J111").Value = 2416
$ws.Range("K111").Value = 7246.799999999999
$ws.Range("L111").Value = 7248
$ws.Range("M111").Value = -4179.799999999999
$ws.Range("N111").Value = -13382
$ws.Range("H132").Value = 1660.2025
$ws.Range("I132").Value = 1499.4266
$ws.Range("J132").Value = 4674.75
$ws.Range("K132").Value = 4498.2798
$ws.Range("L132").Value = 14024.25
$ws.Range("M132").Value = -1968.2798
$ws.Range("N132").Value = -19084.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 889041.8
$ws.Range("I32").Value = 1033980.56
$ws.Range("J32").Value = 19409.223
$ws.Range("K32").Value = 1033980.56
$ws.Range("L32").Value = 19409.223
$ws.Range("M32").Value = -1033693.56
$ws.Range("N32").Value = -19983.223
$ws.Range("H60").Value = 55555
$ws.Range("I60").Value = 55555
$ws.Range("K60").Value = 55555
$ws.Range("M60").Value = -54822
$ws.Range("H132").Value = 2504.5894
$ws.Range("I132").Value = 1417.6136
$ws.Range("J132").Value = 6490.1665
$ws.Range("K132").Value = 4252.8408
$ws.Range("L132").Value = 19470.4995
$ws.Range("M132").Value = -1722.8408
$ws.Range("N132").Value = -24530.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4166.25
$ws.Range("I31").Value = 996.2143
$ws.Range("J31").Value = 8604.299999999999
$ws.Range("K31").Value = 996.2143
$ws.Range("L31").Value = 8604.299999999999
$ws.Range("M31").Value = -701.2143
$ws.Range("N31").Value = -9194.299999999999
$ws.Range("H34").Value = 4166.25
$ws.Range("I34").Value = 996.2143
$ws.Range("J34").Value = 8604.299999999999
$ws.Range("K34").Value = 996.2143
$ws.Range("L34").Value = 8604.299999999999
$ws.Range("M34").Value = -794.2143
$ws.Range("N34").Value = -9008.299999999999
$ws.Range("H62").Value = 2801.25
$ws.Range("I62").Value = 2743.3872
$ws.Range("J62").Value = 3160
$ws.Range("K62").Value = 2743.3872
$ws.Range("L62").Value = 3160
$ws.Range("M62").Value = -2119.3872
$ws.Range("N62").Value = -4408
$ws.Range("H65").Value = 2801.25
$ws.Range("I65").Value = 2743.3872
$ws.Range("J65").Value = 3160
$ws.Range("K65").Value = 13716.936
$ws.Range("L65").Value = 15800
$ws.Range("M65").Value = -10596.936
$ws.Range("N65").Value = -22040
$ws.Range("H99").Value = 1972.2222
$ws.Range("H105").Value = 1169.75
$ws.Range("I105").Value = 899
$ws.Range("J105").Value = 1260
$ws.Range("K105").Value = 899
$ws.Range("L105").Value = 1260
$ws.Range("M105").Value = 848
$ws.Range("N105").Value = -4754
$ws.Range("H126").Value = 1972.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1133.3334
$ws.Range("I57").Value = 500
$ws.Range("K57").Value = 1500
$ws.Range("M57").Value = -941
$ws.Range("H63").Value = 3382.4
$ws.Range("I63").Value = 912
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 2736
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -1987
$ws.Range("N63").Value = -13498
$ws.Range("H66").Value = 3382.4
$ws.Range("I66").Value = 912
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 8208
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -4464
$ws.Range("N66").Value = -43488
$ws.Range("H137").Value = 6179218.5
$ws.Range("I137").Value = 12829392
$ws.Range("J137").Value = 4057.1428
$ws.Range("K137").Value = 38488176
$ws.Range("L137").Value = 12171.4284
$ws.Range("M137").Value = -38483076
$ws.Range("N137").Value = -22371.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 185.06667
$ws.Range("I107").Value = 151.77777
$ws.Range("J107").Value = 235
$ws.Range("K107").Value = 151.77777
$ws.Range("L107").Value = 235
$ws.Range("M107").Value = 1768.22223
$ws.Range("N107").Value = -4075
$ws.Range("H126").Value = 3046.182
$ws.Range("I126").Value = 3000.8
$ws.Range("K126").Value = 9002.400000000001
$ws.Range("M126").Value = -6532.400000000001
$ws.Range("H132").Value = 2683.173
$ws.Range("I132").Value = 2458.578
$ws.Range("J132").Value = 4127
$ws.Range("K132").Value = 7375.734
$ws.Range("L132").Value = 12381
$ws.Range("M132").Value = -4845.734
$ws.Range("N132").Value = -17441

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5594.875
$ws.Range("I7").Value = 5668
$ws.Range("J7").Value = 5551
$ws.Range("K7").Value = 5668
$ws.Range("L7").Value = 5551
$ws.Range("M7").Value = -5556
$ws.Range("N7").Value = -5775
$ws.Range("H58").Value = 7500
$ws.Range("I58").Value = 7500
$ws.Range("K58").Value = 7500
$ws.Range("M58").Value = -7240
$ws.Range("H68").Value = 1582.6666
$ws.Range("I68").Value = 1570.2858
$ws.Range("K68").Value = 1570.2858
$ws.Range("M68").Value = -821.2858000000001
$ws.Range("H71").Value = 1582.6666
$ws.Range("I71").Value = 1570.2858
$ws.Range("K71").Value = 7851.429
$ws.Range("M71").Value = -4107.429
$ws.Range("H126").Value = 5594.875
$ws.Range("I126").Value = 5668
$ws.Range("J126").Value = 5551
$ws.Range("K126").Value = 17004
$ws.Range("L126").Value = 16653
$ws.Range("M126").Value = -14534
$ws.Range("N126").Value = -21593

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1899.1915
$ws.Range("I122").Value = 1922.2778
$ws.Range("J122").Value = 1823.6364
$ws.Range("K122").Value = 5766.8334
$ws.Range("L122").Value = 5470.9092
$ws.Range("M122").Value = -3316.8334
$ws.Range("N122").Value = -10370.9092
$ws.Range("H132").Value = 1399.3934
$ws.Range("I132").Value = 1253.289
$ws.Range("K132").Value = 3759.867
$ws.Range("M132").Value = -1229.867

Write-Host "Applied 168 cell updates across sheets."
